$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new keyword set
$ws.Range("B1").Value = "Processing"
$ws.Range("C1").Value = "AI"
$ws.Range("D1").Value = "Technology"
$ws.Range("E1").Value = "Data"
$ws.Range("F1").Value = "Adjectives"

# Row 2 (2019) data
$ws.Range("B2").Value = 46
$ws.Range("C2").Value = 5242
$ws.Range("D2").Value = 702
$ws.Range("E2").Value = 287
$ws.Range("F2").Value = 28

# Row 3 (2020) data
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 269
$ws.Range("D3").Value = 118
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = 2

# Remove the now-unused last column (G) entirely, shrinking the sheet to A1:F3
$ws.Range("G1:G3").EntireColumn.Delete()
